$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1186.7894
$ws.Range("J17").Value = 1186.7894
$ws.Range("L17").Value = 3560.3682
$ws.Range("N17").Value = -3896.3682

$ws.Range("H62").Value = 8000
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 8000
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -56240

$ws.Range("H92").Value = 1294.4286
$ws.Range("I92").Value = 1348.9375
$ws.Range("K92").Value = 1348.9375
$ws.Range("M92").Value = -100.9375

$ws.Range("H100").Value = 1687.9565
$ws.Range("I100").Value = 1327.875
$ws.Range("J100").Value = 1880
$ws.Range("K100").Value = 1327.875
$ws.Range("L100").Value = 1880
$ws.Range("M100").Value = -786.875
$ws.Range("N100").Value = -2962

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3525.987
$ws.Range("I32").Value = 3187.6943
$ws.Range("K32").Value = 3187.6943
$ws.Range("M32").Value = -2900.6943

$ws.Range("H45").Value = 1210.3846
$ws.Range("I45").Value = 1148.2941
$ws.Range("J45").Value = 1327.6666
$ws.Range("K45").Value = 1148.2941
$ws.Range("L45").Value = 1327.6666
$ws.Range("M45").Value = -771.2941000000001
$ws.Range("N45").Value = -2081.6666

$ws.Range("H97").Value = 536.3158
$ws.Range("I97").Value = 424.375
$ws.Range("J97").Value = 1133.3334
$ws.Range("K97").Value = 424.375
$ws.Range("L97").Value = 1133.3334
$ws.Range("M97").Value = 71.625
$ws.Range("N97").Value = -2125.3334

$ws.Range("H110").Value = 1715.4117
$ws.Range("I110").Value = 1231.7273
$ws.Range("K110").Value = 1231.7273
$ws.Range("M110").Value = 813.2727

$ws.Range("H128").Value = 74614.5
$ws.Range("J128").Value = 74614.5
$ws.Range("L128").Value = 74614.5
$ws.Range("N128").Value = -84574.5

$ws.Range("H132").Value = 1610.4186
$ws.Range("I132").Value = 1471.5416
$ws.Range("K132").Value = 4414.6248
$ws.Range("M132").Value = -1884.6248

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1400.3334
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 1500.5
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 1500.5
$ws.Range("M22").Value = -1027
$ws.Range("N22").Value = -1846.5

$ws.Range("H107").Value = 1760.8667
$ws.Range("I107").Value = 1469.2307
$ws.Range("K107").Value = 1469.2307
$ws.Range("M107").Value = 450.7692999999999

$ws.Range("H134").Value = 5026
$ws.Range("I134").Value = 1492.1923
$ws.Range("J134").Value = 15234.777
$ws.Range("K134").Value = 4476.5769
$ws.Range("L134").Value = 45704.331
$ws.Range("M134").Value = -1941.5769
$ws.Range("N134").Value = -50774.331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45555.223
$ws.Range("J20").Value = 45555.223
$ws.Range("L20").Value = 45555.223
$ws.Range("N20").Value = -46027.223

$ws.Range("H30").Value = 45555.223
$ws.Range("J30").Value = 45555.223
$ws.Range("L30").Value = 45555.223
$ws.Range("N30").Value = -45737.223

$ws.Range("H31").Value = 1380.5
$ws.Range("I31").Value = 1380.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1380.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1085.5
$ws.Range("N31").Value = ""

$ws.Range("H34").Value = 1380.5
$ws.Range("I34").Value = 1380.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1380.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1178.5
$ws.Range("N34").Value = ""

$ws.Range("H114").Value = 27995.4
$ws.Range("I114").Value = 21000
$ws.Range("J114").Value = 29744.25
$ws.Range("K114").Value = 21000
$ws.Range("L114").Value = 29744.25
$ws.Range("M114").Value = -16661
$ws.Range("N114").Value = -38422.25

$ws.Range("H128").Value = 45555.223
$ws.Range("J128").Value = 45555.223
$ws.Range("L128").Value = 45555.223
$ws.Range("N128").Value = -55515.223

$ws.Range("H132").Value = 1210.5
$ws.Range("I132").Value = 763.08826
$ws.Range("K132").Value = 2289.26478
$ws.Range("M132").Value = 240.73522

$ws.Range("H134").Value = 1095.5883
$ws.Range("I134").Value = 972.9
$ws.Range("K134").Value = 2918.7
$ws.Range("M134").Value = -383.6999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 105.71429
$ws.Range("I2").Value = 46.666668
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 280.000008
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = -167.000008
$ws.Range("N2").Value = -1126

$ws.Range("H34").Value = 1625
$ws.Range("I34").Value = 1166.6666
$ws.Range("K34").Value = 3499.9998
$ws.Range("M34").Value = -3415.9998

$ws.Range("H123").Value = 2848.1667
$ws.Range("I123").Value = 2696.6667
$ws.Range("J123").Value = 2898.6667
$ws.Range("K123").Value = 8090.000100000001
$ws.Range("L123").Value = 8696.000100000001
$ws.Range("M123").Value = -5640.000100000001
$ws.Range("N123").Value = -13596.0001

$ws.Range("H140").Value = 34296.484
$ws.Range("I140").Value = 94609.164
$ws.Range("J140").Value = 2829
$ws.Range("K140").Value = 283827.492
$ws.Range("L140").Value = 8487
$ws.Range("M140").Value = -278647.492
$ws.Range("N140").Value = -18847

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1693.439
$ws.Range("I132").Value = 1149.7407
$ws.Range("K132").Value = 3449.2221
$ws.Range("M132").Value = -919.2221

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 999.5
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""

$ws.Range("H84").Value = 999.5
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""

$ws.Range("H100").Value = 859.5
$ws.Range("J100").Value = 567.8
$ws.Range("L100").Value = 1135.6
$ws.Range("N100").Value = -2217.6

$ws.Range("H132").Value = 1009.9787
$ws.Range("I132").Value = 822.81396
$ws.Range("J132").Value = 3022
$ws.Range("K132").Value = 2468.44188
$ws.Range("L132").Value = 9066
$ws.Range("M132").Value = 61.55812000000014
$ws.Range("N132").Value = -14126

$ws.Range("H136").Value = 375.7
$ws.Range("I136").Value = 200.83333
$ws.Range("J136").Value = 1949.5
$ws.Range("K136").Value = 602.49999
$ws.Range("L136").Value = 5848.5
$ws.Range("M136").Value = 1947.50001
$ws.Range("N136").Value = -10948.5
